$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated crypto price/volume values. Values are entered with a
# leading apostrophe to force text (matching source data, which can
# include multi-dot numbers like '60.616.07' that aren't valid numerics),
# then the cell style is reset to Normal so no stray quote-prefix style
# is left behind on the cell.

$ws.Range("D2").Value = "'60.616.07"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -0.50%  "
$ws.Range("E2").Style = "Normal"

$ws.Range("D3").Value = "'2.397.54"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -0.02%  "
$ws.Range("E3").Style = "Normal"

$ws.Range("E4").Value = "'  +0.05%  "
$ws.Range("E4").Style = "Normal"

$ws.Range("D5").Value = "'563.10"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -1.27%  "
$ws.Range("E5").Style = "Normal"

$ws.Range("D6").Value = "'140.77"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +0.39%  "
$ws.Range("E6").Style = "Normal"

$ws.Range("E7").Value = "'  -0.19%  "
$ws.Range("E7").Style = "Normal"

$ws.Range("D8").Value = "'0.534"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  +1.12%  "
$ws.Range("E8").Style = "Normal"

$ws.Range("D9").Value = "'2.402.90"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  +0.21%  "
$ws.Range("E9").Style = "Normal"

$ws.Range("D10").Value = "'0.107"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  -0.86%  "
$ws.Range("E10").Style = "Normal"

$ws.Range("E11").Value = "'  -0.59%  "
$ws.Range("E11").Style = "Normal"

$ws.Range("D12").Value = "'5.16"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  +1.30%  "
$ws.Range("E12").Style = "Normal"

$ws.Range("E13").Value = "'  +0.19%  "
$ws.Range("E13").Style = "Normal"

$ws.Range("D14").Value = "'26.07"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +0.14%  "
$ws.Range("E14").Style = "Normal"

$ws.Range("E15").Value = "'  -0.31%  "
$ws.Range("E15").Style = "Normal"

$ws.Range("E16").Value = "'  -2.04%  "
$ws.Range("E16").Style = "Normal"

$ws.Range("D17").Value = "'60.618.07"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -0.40%  "
$ws.Range("E17").Style = "Normal"

$ws.Range("D18").Value = "'2.406.46"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +0.43%  "
$ws.Range("E18").Style = "Normal"

$ws.Range("D19").Value = "'7.92"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +6.22%  "
$ws.Range("E19").Style = "Normal"

$ws.Range("E20").Value = "'  -0.35%  "
$ws.Range("E20").Style = "Normal"

$ws.Range("D21").Value = "'323.76"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +0.21%  "
$ws.Range("E21").Style = "Normal"

$ws.Range("E22").Value = "'  +0.47%  "
$ws.Range("E22").Style = "Normal"

$ws.Range("D23").Value = "'6.10"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +0.96%  "
$ws.Range("E23").Style = "Normal"

$ws.Range("E24").Value = "'  -0.10%  "
$ws.Range("E24").Style = "Normal"

$ws.Range("D25").Value = "'1.82"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -3.03%  "
$ws.Range("E25").Style = "Normal"

$ws.Range("D26").Value = "'64.99"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  +0.11%  "
$ws.Range("E26").Style = "Normal"

$ws.Range("D27").Value = "'563.68"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -2.93%  "
$ws.Range("E27").Style = "Normal"

$ws.Range("D28").Value = "'8.06"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -4.45%  "
$ws.Range("E28").Style = "Normal"

$ws.Range("D29").Value = "'2.517.98"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -0.26%  "
$ws.Range("E29").Style = "Normal"

$ws.Range("D30").Value = "'0.0₃0934"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  +1.25%  "
$ws.Range("E30").Style = "Normal"

$ws.Range("D31").Value = "'8.08"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  +1.76%  "
$ws.Range("E31").Style = "Normal"

$ws.Range("E32").Value = "'  -1.67%  "
$ws.Range("E32").Style = "Normal"

$ws.Range("D33").Value = "'1.81"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  -1.64%  "
$ws.Range("E33").Style = "Normal"

$ws.Range("E34").Value = "'  -1.44%  "
$ws.Range("E34").Style = "Normal"

$ws.Range("E35").Value = "'  -0.21%  "
$ws.Range("E35").Style = "Normal"

$ws.Range("D36").Value = "'1.45"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  +3.31%  "
$ws.Range("E36").Style = "Normal"

$ws.Range("D37").Value = "'152.27"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +1.48%  "
$ws.Range("E37").Style = "Normal"

$ws.Range("D38").Value = "'0.370"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  +0.76%  "
$ws.Range("E38").Style = "Normal"

$ws.Range("D39").Value = "'4.58"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -1.49%  "
$ws.Range("E39").Style = "Normal"

$ws.Range("D40").Value = "'18.25"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -0.17%  "
$ws.Range("E40").Style = "Normal"

$ws.Range("E41").Value = "'  -0.33%  "
$ws.Range("E41").Style = "Normal"

$ws.Range("E42").Value = "'  -0.05%  "
$ws.Range("E42").Style = "Normal"

$ws.Range("D43").Value = "'41.67"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +1.20%  "
$ws.Range("E43").Style = "Normal"

$ws.Range("D44").Value = "'1.67"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -1.22%  "
$ws.Range("E44").Style = "Normal"

$ws.Range("D45").Value = "'2.48"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +5.67%  "
$ws.Range("E45").Style = "Normal"

$ws.Range("D46").Value = "'0.0₆0283"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -1.52%  "
$ws.Range("E46").Style = "Normal"

$ws.Range("D47").Value = "'141.29"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -0.09%  "
$ws.Range("E47").Style = "Normal"

$ws.Range("D48").Value = "'3.54"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  +0.36%  "
$ws.Range("E48").Style = "Normal"

$ws.Range("E49").Value = "'  -0.10%  "
$ws.Range("E49").Style = "Normal"

$ws.Range("E50").Value = "'  +0.39%  "
$ws.Range("E50").Style = "Normal"

$ws.Range("D51").Value = "'19.29"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -1.45%  "
$ws.Range("E51").Style = "Normal"

